$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "JavaFile TestCase"

# Update header row (row 1)
$ws.Range("A1").Value = "Method Name"
$ws.Range("B1").Value = "Inputs"
$ws.Range("C1").Value = "Expected Output"
$ws.Range("D1").Value = "Expected Status Code"

# New "Expected Status Code" data column (D2:D5) - all 200
$ws.Range("D2").Value = 200
$ws.Range("D3").Value = 200
$ws.Range("D4").Value = 200
$ws.Range("D5").Value = 200

# Header row formatting: vertical center alignment, taller row
$headerRow = $ws.Rows.Item(1)
$headerRow.VerticalAlignment = -4108
$headerRow.RowHeight = 31.5

# Update the selected cell
$ws.Range("C9").Select()
